# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows at the top of the "Choclo" data
# block (rows 502-503), pushing the existing rows 502-509 down to 504-511.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 502 (shifts old rows 502:509 -> 504:511,
# carrying their formatting, including the date-format style on column D).
$ws.Rows("502:503").Insert()

# --- New row 502 : Choclero, Provincia del Elquí -------------------------
$ws.Cells.Item(502, 1).Value  = 8
$ws.Cells.Item(502, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(502, 3).Value  = "Coquimbo"
$ws.Cells.Item(502, 4).Value  = 44595
$ws.Cells.Item(502, 5).Value  = 4
$ws.Cells.Item(502, 6).Value  = 100112024
$ws.Cells.Item(502, 7).Value  = "Choclo"
$ws.Cells.Item(502, 8).Value  = "Choclero"
$ws.Cells.Item(502, 9).Value  = "Primera"
$ws.Cells.Item(502, 10).Value = 20000
$ws.Cells.Item(502, 11).Value = 200
$ws.Cells.Item(502, 12).Value = 250
$ws.Cells.Item(502, 13).Value = 225
$ws.Cells.Item(502, 14).Value = "$/unidad"
$ws.Cells.Item(502, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(502, 16).Value = 225
$ws.Cells.Item(502, 17).Value = 1
$ws.Cells.Item(502, 18).Value = "Hortaliza"

# --- New row 503 : Dulce o Americano, Provincia del Elquí ----------------
$ws.Cells.Item(503, 1).Value  = 8
$ws.Cells.Item(503, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(503, 3).Value  = "Coquimbo"
$ws.Cells.Item(503, 4).Value  = 44595
$ws.Cells.Item(503, 5).Value  = 4
$ws.Cells.Item(503, 6).Value  = 100112024
$ws.Cells.Item(503, 7).Value  = "Choclo"
$ws.Cells.Item(503, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(503, 9).Value  = "Primera"
$ws.Cells.Item(503, 10).Value = 30000
$ws.Cells.Item(503, 11).Value = 150
$ws.Cells.Item(503, 12).Value = 160
$ws.Cells.Item(503, 13).Value = 155
$ws.Cells.Item(503, 14).Value = "$/unidad"
$ws.Cells.Item(503, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(503, 16).Value = 155
$ws.Cells.Item(503, 17).Value = 1
$ws.Cells.Item(503, 18).Value = "Hortaliza"
